$d = $word.ActiveDocument

# --- 1. Collapse the three CORE COMPETENCIES paragraphs into one short line ---
$p6 = $d.Paragraphs.Item(6)
$p6.Range.Text = "Survey Methodology & Research Design " + [char]8226 + " Redistricting & Geospatial Analysis " + [char]8226 + " Data Analysis & Visualization"

# Remove the two now-redundant long paragraphs that followed it
$p7 = $d.Paragraphs.Item(7)
$p9 = $d.Paragraphs.Item(9)
$rng = $d.Range($p7.Range.Start, $p9.Range.Start)
$rng.Delete()

# --- 2. Insert a new TECHNICAL SKILLS section before the closing paragraph ---
$lastBullet = $d.Paragraphs.Item($d.Paragraphs.Count - 1)

$lastBullet.Range.InsertParagraphAfter()
$heading = $d.Paragraphs.Item($lastBullet.Index + 1)
$heading.Style = "Heading2"
$heading.Range.Text = "TECHNICAL SKILLS"

$heading.Range.InsertParagraphAfter()
$skill1 = $d.Paragraphs.Item($heading.Index + 1)
$skill1.Style = "Normal"
$skill1.Range.Text = "SURVEY METHODOLOGY & RESEARCH DESIGN Survey Design and Questionnaire Development for Political and Market Research; Sampling Methodology and Statistical Analysis (R, SPSS, Stata, OSCAR); Random Device Engagement (RDE), Text Message, Web Panel, and Live Telephone Calling; Expert Testimony and Consultation on Research Methodology"

$skill1.Range.InsertParagraphAfter()
$skill2 = $d.Paragraphs.Item($skill1.Index + 1)
$skill2.Style = "Normal"
$skill2.Range.Text = "REDISTRICTING & GEOSPATIAL ANALYSIS Redistricting Software Development and Boundary Estimation Systems; Geospatial Analysis; Choropleths and Hexagonal Grid Maps for Demographic Visualization; Court Case Analysis and Expert Testimony for Redistricting"

$skill2.Range.InsertParagraphAfter()
$skill3 = $d.Paragraphs.Item($skill2.Index + 1)
$skill3.Style = "Normal"
$skill3.Range.Text = "DATA ANALYSIS & VISUALIZATION Advanced Statistical Modeling and Analysis (Regression, Clustering, Segmentation); Data Visualization; Consumer Behavior Analysis and Market Segmentation; Multi-million Dollar Research Project Management"

Write-Output "done"
